$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "68.799.05"
$ws.Cells.Item(2, 5).Value = "  +0.38%  "
$ws.Cells.Item(3, 4).Value = "2.468.72"
$ws.Cells.Item(3, 5).Value = "  +0.36%  "
$ws.Cells.Item(4, 5).Value = "  +0.05%  "
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = "561.42"
$cell.ClearFormats()
$ws.Cells.Item(5, 5).Value = "  +0.17%  "
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = "163.46"
$cell.ClearFormats()
$ws.Cells.Item(6, 5).Value = "  -0.46%  "
$ws.Cells.Item(7, 5).Value = "  +0.13%  "
$ws.Cells.Item(8, 5).Value = "  +1.73%  "
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.157"
$cell.ClearFormats()
$ws.Cells.Item(9, 5).Value = "  +4.10%  "
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.166"
$cell.ClearFormats()
$ws.Cells.Item(10, 5).Value = "  +0.67%  "
$ws.Cells.Item(11, 5).Value = "  -1.91%  "
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.87"
$cell.ClearFormats()
$ws.Cells.Item(12, 5).Value = "  +0.85%  "
$ws.Cells.Item(13, 4).Value = "68.730.48"
$ws.Cells.Item(13, 5).Value = "  +0.47%  "
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0000170"
$cell.ClearFormats()
$ws.Cells.Item(14, 5).Value = "  -0.10%  "
$ws.Cells.Item(15, 5).Value = "  +0.51%  "
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = "10.66"
$cell.ClearFormats()
$ws.Cells.Item(16, 5).Value = "  -3.28%  "
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = "337.11"
$cell.ClearFormats()
$ws.Cells.Item(17, 5).Value = "  -2.15%  "
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = "6.95"
$cell.ClearFormats()
$ws.Cells.Item(18, 5).Value = "  -3.75%  "
$ws.Cells.Item(19, 5).Value = "  -0.12%  "
$ws.Cells.Item(20, 5).Value = "  +0.64%  "
$ws.Cells.Item(21, 5).Value = "  +0.10%  "
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = "66.58"
$cell.ClearFormats()
$ws.Cells.Item(22, 5).Value = "  -1.86%  "
$ws.Cells.Item(23, 5).Value = "  -1.94%  "
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = "8.28"
$cell.ClearFormats()
$ws.Cells.Item(24, 5).Value = "  +0.89%  "
$ws.Cells.Item(25, 4).Value = "0.0₃0826"
$ws.Cells.Item(25, 5).Value = "  -1.61%  "
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.22"
$cell.ClearFormats()
$ws.Cells.Item(26, 5).Value = "  -0.90%  "
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.ClearFormats()
$ws.Cells.Item(27, 5).Value = "  -0.04%  "
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = "430.51"
$cell.ClearFormats()
$ws.Cells.Item(28, 5).Value = "  -0.99%  "
$ws.Cells.Item(29, 5).Value = "  -2.18%  "
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.63"
$cell.ClearFormats()
$ws.Cells.Item(30, 5).Value = "  -3.42%  "
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = "159.65"
$cell.ClearFormats()
$ws.Cells.Item(31, 5).Value = "  +1.07%  "
$ws.Cells.Item(32, 5).Value = "  +0.07%  "
$ws.Cells.Item(33, 5).Value = "  -0.03%  "
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.109"
$cell.ClearFormats()
$ws.Cells.Item(34, 5).Value = "  -1.45%  "
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = "17.85"
$cell.ClearFormats()
$ws.Cells.Item(35, 5).Value = "  -0.38%  "
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.45"
$cell.ClearFormats()
$ws.Cells.Item(36, 5).Value = "  -0.70%  "
$ws.Cells.Item(37, 5).Value = "  -2.25%  "
$ws.Cells.Item(38, 5).Value = "  -3.46%  "
$ws.Cells.Item(39, 5).Value = "  -1.63%  "
$ws.Cells.Item(40, 5).Value = "  -0.80%  "
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.37"
$cell.ClearFormats()
$ws.Cells.Item(41, 5).Value = "  +0.04%  "
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = "130.25"
$cell.ClearFormats()
$ws.Cells.Item(42, 5).Value = "  -3.49%  "
$ws.Cells.Item(43, 5).Value = "  +0.19%  "
$ws.Cells.Item(44, 5).Value = "  -0.24%  "
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.565"
$cell.ClearFormats()
$ws.Cells.Item(45, 5).Value = "  +0.24%  "
$ws.Cells.Item(46, 5).Value = "  +0.76%  "
$ws.Cells.Item(47, 5).Value = "  +0.30%  "
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.40"
$cell.ClearFormats()
$ws.Cells.Item(48, 5).Value = "  -2.55%  "
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = "5.02"
$cell.ClearFormats()
$ws.Cells.Item(49, 5).Value = "  -6.98%  "
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = "16.92"
$cell.ClearFormats()
$ws.Cells.Item(50, 5).Value = "  -4.19%  "
$ws.Cells.Item(51, 4).Value = "0.0₆0208"
$ws.Cells.Item(51, 5).Value = "  +0.45%  "
